# The commit reorganises the test-data files and, for this workbook,
# duplicates the "status"/"Pass" column (F) into a brand new column G on
# the active "RegTestData" sheet (dimension grows from A1:F4 to A1:G4,
# sharedStrings gains 4 more references to the existing "status"/"Pass"
# strings, and a duplicate cellXfs entry is produced for the new header
# cell's fill style).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegTestData")

# Copy the header + 3 data rows of column F ("status" / "Pass" x3) into
# the new column G, reusing the same shared-string values as F.
$ws.Range("G1:G4").Value2 = $ws.Range("F1:F4").Value2

# F1 carries a shaded fill (style index 3). Re-applying that same fill to
# G1 registers a fresh (duplicate) style entry for the new header cell,
# matching the extra <xf> row the original edit introduced.
$ws.Range("G1").Interior.ColorIndex = $ws.Range("F1").Interior.ColorIndex
